$wb = $excel.ActiveWorkbook

# "Repayment schedule" sheet (3rd tab) gains a new blank column before column N ("Late"),
# inheriting the width/format of the column to its left (mirrors Excel's native
# Insert Column behaviour).
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Columns("N").Insert()

# Update the active selection on that sheet and make it the active tab.
$wsSchedule.Range("S6").Select()

# "Summary" sheet was previously the active tab; it no longer is.
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Select()

$wsSchedule.Select()
